$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "69.750.07"
Set-TextValue "E2" "  -1.29%  "
Set-TextValue "D3" "3.566.66"
Set-TextValue "E3" "  -2.32%  "
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "574.65"
Set-TextValue "E5" "  -3.43%  "
Set-TextValue "D6" "186.81"
Set-TextValue "E6" "  -3.81%  "
Set-TextValue "D7" "3.562.03"
Set-TextValue "E7" "  -2.20%  "
Set-TextValue "D8" "0.620"
Set-TextValue "E8" "  -4.05%  "
Set-TextValue "E9" "  +0.03%  "
Set-TextValue "D10" "0.182"
Set-TextValue "E10" "  -0.82%  "
Set-TextValue "E11" "  -4.20%  "
Set-TextValue "D12" "54.69"
Set-TextValue "E12" "  -5.98%  "
Set-TextValue "D13" "0.0000299"
Set-TextValue "E13" "  +2.10%  "
Set-TextValue "D14" "9.51"
Set-TextValue "E14" "  -4.39%  "
Set-TextValue "D15" "4.140.22"
Set-TextValue "E15" "  -2.08%  "
Set-TextValue "D16" "19.54"
Set-TextValue "E16" "  -2.65%  "
Set-TextValue "D17" "3.567.07"
Set-TextValue "E17" "  -2.14%  "
Set-TextValue "D18" "69.732.75"
Set-TextValue "E18" "  -1.30%  "
Set-TextValue "D19" "12.52"
Set-TextValue "E19" "  -1.86%  "
Set-TextValue "E20" "  -0.65%  "
Set-TextValue "E21" "  -3.66%  "
Set-TextValue "D22" "485.51"
Set-TextValue "E22" "  -0.72%  "
Set-TextValue "D23" "19.31"
Set-TextValue "E23" "  +0.81%  "
Set-TextValue "E24" "  -7.59%  "
Set-TextValue "D25" "4.37"
Set-TextValue "E25" "  -3.89%  "
Set-TextValue "D26" "95.05"
Set-TextValue "E26" "  +4.18%  "
Set-TextValue "D27" "11.37"
Set-TextValue "E27" "  -1.15%  "
Set-TextValue "E28" "  -6.84%  "
Set-TextValue "E29" "  -3.36%  "
Set-TextValue "E30" "  -3.56%  "
Set-TextValue "D32" "66.71"
Set-TextValue "E32" "  +1.29%  "
Set-TextValue "D33" "11.99"
Set-TextValue "E33" "  -2.10%  "
Set-TextValue "E34" "  -6.18%  "
Set-TextValue "D35" "566.67"
Set-TextValue "E35" "  -9.39%  "
Set-TextValue "D36" "3.16"
Set-TextValue "E36" "  +12.01%  "
Set-TextValue "D37" "38.48"
Set-TextValue "E37" "  -5.03%  "
Set-TextValue "D38" "1.00"
Set-TextValue "E38" "  -0.03%  "
Set-TextValue "D39" "0.396"
Set-TextValue "E39" "  -3.85%  "
Set-TextValue "D40" "0.0₃0792"
Set-TextValue "E40" "  -4.04%  "
Set-TextValue "D41" "3.48"
Set-TextValue "E41" "  -2.84%  "
Set-TextValue "D42" "3.17"
Set-TextValue "E42" "  +4.18%  "
Set-TextValue "D43" "0.135"
Set-TextValue "E43" "  -8.65%  "
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "3.230.26"
Set-TextValue "E44" "  -2.18%  "
Set-TextValue "B45" "ThetaToken"
Set-TextValue "C45" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D45" "2.99"
Set-TextValue "E45" "  -5.22%  "
Set-TextValue "E46" "  -4.32%  "
Set-TextValue "B47" "ApeXProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "3.40"
Set-TextValue "E47" "  +2.77%  "
Set-TextValue "B48" "THORChain"
Set-TextValue "C48" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D48" "9.61"
Set-TextValue "E48" "  +3.05%  "
Set-TextValue "D49" "0.135"
Set-TextValue "E49" "  -2.50%  "
Set-TextValue "E50" "  +0.24%  "
Set-TextValue "D51" "3.18"
Set-TextValue "E51" "  -3.86%  "
